$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = "0.044`n (0.053)"
$ws.Range("D2").Value = "0.011`n (0.014)"
$ws.Range("F2").Value = "-0.011`n (0.011)"
$ws.Range("G2").Value = "0.019`n (0.023)"

# Row 3
$ws.Range("B3").Value = "-0.370***`n (0.115)"
$ws.Range("C3").Value = "-0.589**`n (0.226)"
$ws.Range("D3").Value = "0.468***`n (0.082)"
$ws.Range("E3").Value = "0.236`n (0.147)"
$ws.Range("F3").Value = "0.343***`n (0.061)"
$ws.Range("G3").Value = "0.141`n (0.099)"

# Row 4
$ws.Range("B4").Value = "0.308*`n (0.181)"
$ws.Range("D4").Value = "0.071`n (0.126)"
$ws.Range("F4").Value = "0.151`n (0.100)"

# Row 5 (numeric R-squared values)
$ws.Range("B5").Value = 0.2857796200381378
$ws.Range("C5").Value = 0.1516449335735305
$ws.Range("D5").Value = 0.4886088471351145
$ws.Range("E5").Value = 0.06520864925664716
$ws.Range("F5").Value = 0.3805361884847204
$ws.Range("G5").Value = 0.0551647097794542
